$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '29.239.57'
$ws.Range('E2').Value = '  -0.33%  '
Set-TextValue $ws.Range('D3') '1.838.09'
$ws.Range('E3').Value = '  -0.40%  '
Set-TextValue $ws.Range('D4') '0.9977'
$ws.Range('E4').Value = '  -0.26%  '
Set-TextValue $ws.Range('D5') '240.81'
$ws.Range('E5').Value = '  -1.42%  '
Set-TextValue $ws.Range('D6') '0.6702'
$ws.Range('E6').Value = '  -2.44%  '
Set-TextValue $ws.Range('D7') '0.9989'
$ws.Range('E7').Value = '  -0.17%  '
Set-TextValue $ws.Range('D8') '0.07413'
$ws.Range('E8').Value = '  -1.73%  '
Set-TextValue $ws.Range('D9') '0.2959'
$ws.Range('E9').Value = '  -2.38%  '
Set-TextValue $ws.Range('D10') '22.95'
$ws.Range('E10').Value = '  -1.58%  '
Set-TextValue $ws.Range('D11') '0.07727'
$ws.Range('E11').Value = '  +0.77%  '
Set-TextValue $ws.Range('D12') '1.860.63'
$ws.Range('E12').Value = '  +0.82%  '
Set-TextValue $ws.Range('D13') '5.031'
$ws.Range('E13').Value = '  -1.23%  '
Set-TextValue $ws.Range('D14') '0.6798'
$ws.Range('E14').Value = '  -1.05%  '
Set-TextValue $ws.Range('D15') '86.52'
$ws.Range('E15').Value = '  -3.16%  '
Set-TextValue $ws.Range('D16') '6.206'
$ws.Range('E16').Value = '  -1.43%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue $ws.Range('D17') '29.268.00'
$ws.Range('E17').Value = '  -0.25%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws.Range('D18') '0.000008263'
$ws.Range('E18').Value = '  +0.75%  '
Set-TextValue $ws.Range('D19') '229.46'
$ws.Range('E19').Value = '  -2.13%  '
Set-TextValue $ws.Range('D20') '12.56'
$ws.Range('E20').Value = '  -0.33%  '
Set-TextValue $ws.Range('D21') '0.9993'
$ws.Range('E21').Value = '  -0.14%  '
Set-TextValue $ws.Range('D22') '7.309'
$ws.Range('E22').Value = '  -3.71%  '
Set-TextValue $ws.Range('D23') '0.9986'
$ws.Range('E23').Value = '  -0.21%  '
Set-TextValue $ws.Range('D24') '160.32'
$ws.Range('E24').Value = '  +0.32%  '
Set-TextValue $ws.Range('D25') '8.731'
$ws.Range('E25').Value = '  -1.49%  '
Set-TextValue $ws.Range('D26') '0.1413'
$ws.Range('E26').Value = '  -3.23%  '
Set-TextValue $ws.Range('D27') '18.03'
$ws.Range('E27').Value = '  -0.18%  '
$ws.Range('E28').Value = '  -1.03%  '
Set-TextValue $ws.Range('D29') '4.215'
$ws.Range('E29').Value = '  -0.43%  '
Set-TextValue $ws.Range('D30') '4.092'
$ws.Range('E30').Value = '  -0.72%  '
$ws.Range('E31').Value = '  -0.77%  '
Set-TextValue $ws.Range('D32') '0.05356'
$ws.Range('E32').Value = '  +3.32%  '
Set-TextValue $ws.Range('D33') '1.875'
$ws.Range('E33').Value = '  +0.59%  '
Set-TextValue $ws.Range('D34') '0.7524'
$ws.Range('E34').Value = '  -2.17%  '
Set-TextValue $ws.Range('D35') '1.141'
$ws.Range('E35').Value = '  +0.03%  '
$ws.Range('E36').Value = '  +0.07%  '
Set-TextValue $ws.Range('D37') '1.334.04'
$ws.Range('E37').Value = '  +2.52%  '
Set-TextValue $ws.Range('D38') '0.01803'
$ws.Range('E38').Value = '  -2.62%  '
Set-TextValue $ws.Range('D39') '2.733'
$ws.Range('E39').Value = '  +1.05%  '
Set-TextValue $ws.Range('D40') '0.9242'
$ws.Range('E40').Value = '  -2.11%  '
Set-TextValue $ws.Range('D41') '5.981'
$ws.Range('E41').Value = '  +4.37%  '
Set-TextValue $ws.Range('D42') '0.08366'
$ws.Range('E42').Value = '  +20.92%  '
Set-TextValue $ws.Range('D43') '0.9990'
Set-TextValue $ws.Range('D44') '103.32'
$ws.Range('E44').Value = '  -1.94%  '
$ws.Range('B45').Value = 'RocketPoolETH'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextValue $ws.Range('D45') '2.009.89'
$ws.Range('E45').Value = '  +0.83%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws.Range('D46') '0.00000000124'
$ws.Range('E46').Value = '  +1.00%  '
$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws.Range('D47') '0.5170'
$ws.Range('E47').Value = '  -0.84%  '
Set-TextValue $ws.Range('D48') '64.16'
$ws.Range('E48').Value = '  +1.95%  '
Set-TextValue $ws.Range('D49') '1.766'
$ws.Range('E49').Value = '  -0.15%  '
Set-TextValue $ws.Range('D50') '9.265'
$ws.Range('E50').Value = '  -4.08%  '
Set-TextValue $ws.Range('D51') '0.05955'
$ws.Range('E51').Value = '  +0.38%  '
